$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the cell to keep a text/string type (matching the source data,
    # which stores prices as inline strings like "1.006" or "27.387.28")
    # instead of letting Excel auto-coerce number-looking text into a
    # numeric cell. We flip the number format to Text only for the
    # duration of the write, then restore the cell's original style so the
    # stylesheet/format stays untouched.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row -> @(D, E) price / volume updates (column B/C only change for the
# rows 40-41 swap, handled separately below).
$updates = [ordered]@{
    2  = @("27.387.28", "  +0.81%  ")
    3  = @("1.789.42", "  +1.77%  ")
    4  = @("1.006", "  +1.01%  ")
    5  = @("337.53", "  +0.35%  ")
    6  = @("1.003", "  +0.65%  ")
    7  = @("0.3796", "  +1.30%  ")
    8  = @("0.3443", "  +0.75%  ")
    9  = @("48.69", "  -0.17%  ")
    10 = @("1.195", "  -0.22%  ")
    11 = @("0.07481", "  -1.29%  ")
    12 = @("1.004", "  +0.90%  ")
    13 = @("21.85", "  +5.77%  ")
    14 = @("6.452", "  +0.45%  ")
    15 = @("1.792.85", "  +2.11%  ")
    16 = @("7.054", "  -0.65%  ")
    17 = @("0.00001099", "  +0.31%  ")
    18 = @("0.06663", "  -1.10%  ")
    19 = @("84.60", "  +1.45%  ")
    20 = @("1.002", "  +0.60%  ")
    21 = @("6.521", "  +3.85%  ")
    22 = @("17.29", "  +1.93%  ")
    23 = @("27.385.80", "  +1.10%  ")
    24 = @("12.51", "  -3.52%  ")
    25 = @("2.432", "  -0.68%  ")
    26 = @("2.552", "  +4.02%  ")
    27 = @("1.490", "  -0.44%  ")
    28 = @("21.47", "  +8.47%  ")
    29 = @("153.53", "  +0.67%  ")
    30 = @("1.993.59", "  +2.22%  ")
    31 = @("133.28", "  +0.05%  ")
    32 = @("4.066", "  -1.29%  ")
    33 = @("6.089", "  -0.05%  ")
    34 = @("0.08673", "  +0.17%  ")
    35 = @("13.16", "  +1.05%  ")
    36 = @("1.651", "  -2.52%  ")
    37 = @("5.454", "  -0.85%  ")
    38 = @("0.6893", "  +8.05%  ")
    39 = @("0.06370", "  +0.20%  ")
    42 = @("0.02340", "  -1.10%  ")
    43 = @("1.261", "  +2.64%  ")
    44 = @("14.39", "  -0.26%  ")
    46 = @("0.6434", "  +2.05%  ")
    47 = @("3.867", "  +0.03%  ")
    48 = @("2.127", "  +1.38%  ")
    49 = @("129.66", "  -0.69%  ")
    50 = @("0.07193", "  -0.67%  ")
    51 = @("79.40", "  +0.50%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    Set-TextValue $ws.Range("D$row") $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

# Row 45 only has the Volume(1h) value change; Price stays the same.
$ws.Range("E45").Value = "  +0.52%  "

# Rows 40 and 41 swap coin order: Algorand moves up to row 40, FraxShare
# moves down to row 41 (rank/column A is unaffected - only B/C/D/E move).
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D40") "0.2204"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D41") "8.822"
$ws.Range("E41").Value = "  +2.43%  "
